$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

Set-TextValue "D2" "24.637.62"
Set-TextValue "E2" "  +3.62%  "
Set-TextValue "D3" "1.697.52"
Set-TextValue "E3" "  +2.33%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  +0.19%  "
Set-TextValue "D5" "314.32"
Set-TextValue "E5" "  +2.37%  "
Set-TextValue "E6" "  +0.06%  "
Set-TextValue "D7" "0.3944"
Set-TextValue "E7" "  +1.51%  "
Set-TextValue "D8" "58.52"
Set-TextValue "E8" "  +20.95%  "
Set-TextValue "D9" "0.4044"
Set-TextValue "E9" "  +2.83%  "
Set-TextValue "D10" "1.526"
Set-TextValue "E10" "  +7.08%  "
Set-TextValue "D11" "1.001"
Set-TextValue "E11" "  +0.23%  "
Set-TextValue "D12" "0.08806"
Set-TextValue "E12" "  +2.18%  "
Set-TextValue "D13" "7.238"
Set-TextValue "E13" "  +12.03%  "
Set-TextValue "D14" "23.36"
Set-TextValue "E14" "  +3.21%  "
Set-TextValue "D15" "0.00001333"
Set-TextValue "E15" "  +3.14%  "
Set-TextValue "D16" "7.598"
Set-TextValue "E16" "  +6.06%  "
Set-TextValue "D17" "1.704.40"
Set-TextValue "E17" "  +2.75%  "
Set-TextValue "D18" "100.84"
Set-TextValue "E18" "  +0.65%  "
Set-TextValue "D19" "0.07086"
Set-TextValue "E19" "  +4.69%  "
Set-TextValue "D20" "19.64"
Set-TextValue "E20" "  +3.45%  "
Set-TextValue "D21" "6.761"
Set-TextValue "E21" "  +2.33%  "
Set-TextValue "E22" "  +0.12%  "
Set-TextValue "D23" "14.19"
Set-TextValue "E23" "  +3.48%  "
Set-TextValue "D24" "24.627.19"
Set-TextValue "E24" "  +3.64%  "
Set-TextValue "D25" "2.990"
Set-TextValue "E25" "  +10.45%  "
Set-TextValue "D26" "2.317"
Set-TextValue "E26" "  +0.71%  "
Set-TextValue "E27" "  +3.62%  "
Set-TextValue "D28" "160.30"
Set-TextValue "E28" "  +2.69%  "
Set-TextValue "D29" "5.207"
Set-TextValue "E29" "  +1.92%  "
Set-TextValue "D30" "134.38"
Set-TextValue "E30" "  +3.74%  "
Set-TextValue "D31" "7.370"
Set-TextValue "E31" "  +28.37%  "
Set-TextValue "D32" "1.116"
Set-TextValue "E32" "  -1.99%  "
Set-TextValue "D33" "1.890.79"
Set-TextValue "E33" "  +2.81%  "
Set-TextValue "D34" "7.462"
Set-TextValue "E34" "  +21.91%  "
Set-TextValue "D35" "0.08623"
Set-TextValue "E35" "  +0.64%  "
Set-TextValue "D36" "1.990"
Set-TextValue "E36" "  +5.90%  "
Set-TextValue "D37" "11.13"
Set-TextValue "E37" "  +8.92%  "
Set-TextValue "D38" "0.2758"
Set-TextValue "E38" "  +5.57%  "
Set-TextValue "D39" "14.82"
Set-TextValue "E39" "  -0.50%  "
Set-TextValue "D40" "0.02773"
Set-TextValue "E40" "  +10.65%  "
Set-TextValue "D41" "0.09087"
Set-TextValue "E41" "  +3.55%  "
Set-TextValue "D42" "1.473"
Set-TextValue "E42" "  +2.34%  "
Set-TextValue "D43" "0.7762"
Set-TextValue "E43" "  +4.88%  "
Set-TextValue "D44" "0.7273"
Set-TextValue "E44" "  +4.77%  "
Set-TextValue "D45" "15.53"
Set-TextValue "E45" "  +5.45%  "
Set-TextValue "D46" "2.502"
Set-TextValue "E46" "  +6.28%  "
Set-TextValue "E47" "  +3.81%  "
Set-TextValue "E48" "  +0.06%  "
Set-TextValue "D49" "141.54"
Set-TextValue "E49" "  +1.04%  "
Set-TextValue "D50" "1.288"
Set-TextValue "E50" "  +14.25%  "
Set-TextValue "B51" "Cronos"
Set-TextValue "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.08017"
Set-TextValue "E51" "  +3.94%  "
